# Refresh the team transition-probability matrix (LIU_A) after simulating
# additional games: each row B:S is a distribution of outcome probabilities
# recomputed from updated (larger) game counts per starting state (A..Q row
# headers in column A). Only the cells whose underlying counts changed are
# rewritten; rows with no recorded transitions (5, 14) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2206405693950178
$ws.Cells.Item(2, 3).Value = 0.498220640569395
$ws.Cells.Item(2, 10).Value = 0.01779359430604982
$ws.Cells.Item(2, 16).Value = 0.1672597864768683
$ws.Cells.Item(2, 19).Value = 0.09608540925266904
$ws.Cells.Item(3, 2).Value = 0.01342281879194631
$ws.Cells.Item(3, 3).Value = 0.03355704697986577
$ws.Cells.Item(3, 10).Value = 0.03355704697986577
$ws.Cells.Item(3, 16).Value = 0.6711409395973155
$ws.Cells.Item(3, 19).Value = 0.2483221476510067
$ws.Cells.Item(4, 10).Value = 0.02777777777777778
$ws.Cells.Item(4, 16).Value = 0.6666666666666666
$ws.Cells.Item(4, 19).Value = 0.3055555555555556
$ws.Cells.Item(6, 2).Value = 0.06122448979591837
$ws.Cells.Item(6, 4).Value = 0.01530612244897959
$ws.Cells.Item(6, 6).Value = 0.07142857142857142
$ws.Cells.Item(6, 10).Value = 0.1479591836734694
$ws.Cells.Item(6, 15).Value = 0.01020408163265306
$ws.Cells.Item(6, 17).Value = 0.1887755102040816
$ws.Cells.Item(6, 18).Value = 0.05612244897959184
$ws.Cells.Item(6, 19).Value = 0.4489795918367347
$ws.Cells.Item(7, 2).Value = 0.07111111111111111
$ws.Cells.Item(7, 4).Value = 0.02222222222222222
$ws.Cells.Item(7, 6).Value = 0.01777777777777778
$ws.Cells.Item(7, 10).Value = 0.1288888888888889
$ws.Cells.Item(7, 15).Value = 0.008888888888888889
$ws.Cells.Item(7, 17).Value = 0.1555555555555556
$ws.Cells.Item(7, 18).Value = 0.08
$ws.Cells.Item(7, 19).Value = 0.5155555555555555
$ws.Cells.Item(8, 2).Value = 0.0954356846473029
$ws.Cells.Item(8, 4).Value = 0.01037344398340249
$ws.Cells.Item(8, 6).Value = 0.05186721991701245
$ws.Cells.Item(8, 10).Value = 0.09336099585062241
$ws.Cells.Item(8, 15).Value = 0.01452282157676349
$ws.Cells.Item(8, 17).Value = 0.2074688796680498
$ws.Cells.Item(8, 18).Value = 0.06846473029045644
$ws.Cells.Item(8, 19).Value = 0.45850622406639
$ws.Cells.Item(9, 2).Value = 0.08227848101265822
$ws.Cells.Item(9, 4).Value = 0.0189873417721519
$ws.Cells.Item(9, 6).Value = 0.0759493670886076
$ws.Cells.Item(9, 10).Value = 0.06962025316455696
$ws.Cells.Item(9, 15).Value = 0.0189873417721519
$ws.Cells.Item(9, 17).Value = 0.1582278481012658
$ws.Cells.Item(9, 18).Value = 0.06329113924050633
$ws.Cells.Item(9, 19).Value = 0.5126582278481012
$ws.Cells.Item(10, 2).Value = 0.1092943201376936
$ws.Cells.Item(10, 4).Value = 0.0189328743545611
$ws.Cells.Item(10, 6).Value = 0.0783132530120482
$ws.Cells.Item(10, 10).Value = 0.09896729776247848
$ws.Cells.Item(10, 15).Value = 0.01376936316695353
$ws.Cells.Item(10, 17).Value = 0.2091222030981067
$ws.Cells.Item(10, 18).Value = 0.06368330464716007
$ws.Cells.Item(10, 19).Value = 0.4079173838209983
$ws.Cells.Item(11, 6).Value = 0.005333333333333333
$ws.Cells.Item(11, 7).Value = 0.1493333333333333
$ws.Cells.Item(11, 10).Value = 0.09333333333333334
$ws.Cells.Item(11, 11).Value = 0.2186666666666667
$ws.Cells.Item(11, 12).Value = 0.5093333333333333
$ws.Cells.Item(11, 19).Value = 0.024
$ws.Cells.Item(12, 7).Value = 0.755
$ws.Cells.Item(12, 10).Value = 0.17
$ws.Cells.Item(12, 11).Value = 0.005
$ws.Cells.Item(12, 12).Value = 0.025
$ws.Cells.Item(12, 19).Value = 0.045
$ws.Cells.Item(13, 7).Value = 0.5681818181818182
$ws.Cells.Item(13, 10).Value = 0.3863636363636364
$ws.Cells.Item(13, 19).Value = 0.04545454545454546
$ws.Cells.Item(15, 6).Value = 0.01123595505617977
$ws.Cells.Item(15, 8).Value = 0.1348314606741573
$ws.Cells.Item(15, 9).Value = 0.0898876404494382
$ws.Cells.Item(15, 10).Value = 0.3707865168539326
$ws.Cells.Item(15, 11).Value = 0.09550561797752809
$ws.Cells.Item(15, 13).Value = 0.01685393258426966
$ws.Cells.Item(15, 15).Value = 0.01123595505617977
$ws.Cells.Item(15, 19).Value = 0.2696629213483146
$ws.Cells.Item(16, 6).Value = 0.01183431952662722
$ws.Cells.Item(16, 8).Value = 0.242603550295858
$ws.Cells.Item(16, 9).Value = 0.05325443786982249
$ws.Cells.Item(16, 10).Value = 0.3964497041420119
$ws.Cells.Item(16, 11).Value = 0.1005917159763314
$ws.Cells.Item(16, 13).Value = 0.03550295857988166
$ws.Cells.Item(16, 15).Value = 0.05325443786982249
$ws.Cells.Item(16, 19).Value = 0.106508875739645
$ws.Cells.Item(17, 6).Value = 0.02073732718894009
$ws.Cells.Item(17, 8).Value = 0.2027649769585254
$ws.Cells.Item(17, 9).Value = 0.08755760368663594
$ws.Cells.Item(17, 10).Value = 0.3940092165898618
$ws.Cells.Item(17, 11).Value = 0.119815668202765
$ws.Cells.Item(17, 13).Value = 0.0184331797235023
$ws.Cells.Item(17, 15).Value = 0.04377880184331797
$ws.Cells.Item(17, 19).Value = 0.1129032258064516
$ws.Cells.Item(18, 6).Value = 0.0136986301369863
$ws.Cells.Item(18, 8).Value = 0.1917808219178082
$ws.Cells.Item(18, 9).Value = 0.1027397260273973
$ws.Cells.Item(18, 10).Value = 0.410958904109589
$ws.Cells.Item(18, 11).Value = 0.1027397260273973
$ws.Cells.Item(18, 13).Value = 0.00684931506849315
$ws.Cells.Item(18, 15).Value = 0.0684931506849315
$ws.Cells.Item(18, 19).Value = 0.1027397260273973
$ws.Cells.Item(19, 6).Value = 0.008227374719521317
$ws.Cells.Item(19, 8).Value = 0.2318623784592371
$ws.Cells.Item(19, 9).Value = 0.05833956619296934
$ws.Cells.Item(19, 10).Value = 0.3649962602842184
$ws.Cells.Item(19, 11).Value = 0.1413612565445026
$ws.Cells.Item(19, 13).Value = 0.02094240837696335
$ws.Cells.Item(19, 14).Value = 0.0007479431563201197
$ws.Cells.Item(19, 15).Value = 0.06507105459985041
$ws.Cells.Item(19, 19).Value = 0.1084517576664173
